# Incremento (x100) nos valores da tabela de ranking - colunas E e F (linhas 2-7)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 7; $row++) {
    foreach ($col in @("E", "F")) {
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $cell.Value2 * 100
    }
}
